# feat: citta combination citta panel
#
# Adds two new "thina"/"middha" combination rows for the "lobha4" citta
# group and two more for the "lobha8" citta group on the "Combination"
# sheet, mirroring the combinations that already exist for "dosa2".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Combination")

# --- Insert the lobha4 combination rows (after the existing "lobha4 / mana" row, row 5) ---
$ws.Rows.Item(6).Resize(2).EntireRow.Insert()

$ws.Cells.Item(6, 1).Value = "lobha4"
$ws.Cells.Item(6, 2).Value = "thina"
$ws.Cells.Item(6, 3).Value = "middha"

$ws.Cells.Item(7, 1).Value = "lobha4"
$ws.Cells.Item(7, 2).Value = "mana"
$ws.Cells.Item(7, 3).Value = "thina"
$ws.Cells.Item(7, 4).Value = "middha"

# --- Insert the lobha8 combination rows (after the existing "lobha8 / mana" row, now row 11) ---
$ws.Rows.Item(12).Resize(2).EntireRow.Insert()

$ws.Cells.Item(12, 1).Value = "lobha8"
$ws.Cells.Item(12, 2).Value = "thina"
$ws.Cells.Item(12, 3).Value = "middha"

$ws.Cells.Item(13, 1).Value = "lobha8"
$ws.Cells.Item(13, 2).Value = "mana"
$ws.Cells.Item(13, 3).Value = "thina"
$ws.Cells.Item(13, 4).Value = "middha"

# Move the selection to reflect where the author ended up working
$ws.Range("E11").Select()
